$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) text (or $null if unchanged), new Volume(1h) (E) text (or $null if unchanged)
# Price cells are forced to text first (NumberFormat "@") so numeric-looking strings
# like "0.999" are not silently coerced into Excel numbers, then the style is put back
# to "Normal" so no stray formatting diff is introduced.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "70.311.79"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.564.62"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.10%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "588.24"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "186.77"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.30%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.552.22"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -1.32%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.620"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("E9").Value = "  +0.00%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.200"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +6.32%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.647"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.46%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "54.50"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.35%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000307"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -1.65%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "9.52"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.27%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.119.86"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "19.46"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.66%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "70.139.47"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.558.19"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.30%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.49"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("E20").Value = "  -1.20%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "538.89"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +9.60%  "

$ws.Range("E22").Value = "  -2.19%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "18.11"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -5.96%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "4.65"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +6.37%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "4.90"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.65%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "95.54"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "11.26"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.54%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "3.00"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.52%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.15"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.79%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "32.23"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.29%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.35"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -3.03%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "12.53"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +3.15%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "65.14"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.05%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.93%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "551.28"
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "3.23"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +6.52%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.417"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +5.34%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "38.64"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.53%  "

$ws.Range("E39").Value = "  -0.10%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0769"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -5.40%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.41"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -3.31%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.134"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.66%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "3.359.97"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +4.37%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.10"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -7.74%  "

$ws.Range("E45").Value = "  +6.40%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.98"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.22%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.0444"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "9.22"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -5.92%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.137"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "137.07"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.16%  "
